$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TrancheBounds")

# Add a new row 4 with the "Mean Gas Rate" label and values
$ws.Range("A4").Value = "Mean Gas Rate"
$ws.Range("B4").Value = 0.0033141656254018852
$ws.Range("C4").Value = 0.023229257753825928
$ws.Range("D4").Value = 0.060885670409754518
$ws.Range("E4").Value = 0.12176465229150621
$ws.Range("F4").Value = 0.27712243647722223
$ws.Range("G4").Value = 0.58810554596876741
$ws.Range("H4").Value = 1.787091075755936
$ws.Range("I4").Value = 5.7882521336484194
$ws.Range("J4").Value = 20.990273004878844
$ws.Range("K4").Value = 161.06811496813148

# Size column A to match the bestFit width seen in the diff ("Mean Gas Rate" is
# now the widest entry in column A)
$ws.Columns.Item(1).ColumnWidth = 12.5

# Move the selection like the diff shows (activeCell B6)
$ws.Range("B6").Select() | Out-Null
